$wb = $excel.ActiveWorkbook

# zh-cn sheet: update handoff/handback datetimes for the last (1622dd2c...) row
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-02-15 04:18:58"
$wsZhCn.Range("G5").Value = "2016-02-15 04:19:44"

# de-de sheet: update handoff/handback datetimes for the last (1622dd2c...) row
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-02-15 04:19:12"
$wsDeDe.Range("G5").Value = "2016-02-15 04:20:10"
